$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2) Update the hidden CEIC metadata comment attached to A1 (its text is an
#    opaque base64-encoded, gzip-compressed blob - replace verbatim with the
#    new snapshot from the authored change).
$newCommentText = "Kh0AAB+LCAAAAAAAAAOlWVtvG8cV/isLPrVAyd0lJVsWRhvwJocoKQokVVl5KYa7I2mq5S67MyuJbymQIEWaoigKp0iv6FOKAnWNNgFSu5f/Eliy+9S/0DOXvZFUzFUNwdo55zszZ86c24zQO9cz37gkEaNhsFexa1bFIIEbejQ426vE/LRqP6i846DutUv8QxzhGeEANkAqYLvXjO5Vzjmf75rm1dVV7apRC6Mzs25Ztvlk0B+752SGqzRgHAcuqaRS3tulKg5qe7MB4djDHCvJvUpv3Ku1CXU7QBvgAJ+RqNaKGQ0IY92AU04JE5IRwZy0O4PvqY059dqDmo3MFXqGbMXU9xSugFR0jYNlyYTOiFO37J2qtVNtWBO7vmvbu41G7cHDrfcSwRSI+pjxMYkuqSsJY45ncylu7TQsu27bjQYy14JgrswADhr63ohcUka8NvF9Vsoipj7Apsth1+WMaSEzJ6snur8KjyM8P59Q7pNyaowGLWMWaF2ySRy0H0bEBfvdS6UDcjWMtFkn8z5wJ+c04osOXpSe64iRaDgXRion6qBOGPCmTyJ+NIezJh64AjAcHsUEmXcwM6EOZS580yAmnnOKfZYXKjDRcRhdsDl2yQHEsSnmuAr8EHvgcJwyTt1s0RUGOozCOcwIi7dC39uHWTV4DSOduReAicWyrTC8yLRbx0TyVOX5wpnOME/gK3Q0Pg+vhoG/GMdT5kZ0SrxOK0Gv5SERkFq6HTMezkCLjIQULUcZDMwF/IMYXOagDnHpDPuHPtiRORC9RQJqxjw8pbwd+vEsYIlaS1R0DJuakOt0k+kYDeF8A2H3MOgFCV5Zei2rKDAKr9I1VxnSDjlyk7nJia8ylsEdoCUnuMqRhyJ2uU99KBH548hRi44xPieEr/UKxUEiG+6LouO0FgfxbAoRNoUwu5SrMmRmfASuCu4OejkWFJKq/JlY1q78AT1SNuoG3t24hIlgudxajg28JRKCPfktHwcXQD2m/PygmexlDQcpC9yJX+UhCN65jxeSnFopT0O9wPVjj6ic0AtOpYsK3dSh3slGK6Q+hLmDcLCYLOaQmhnd5fCxV4Fivct4BO1AxXHDOODRQiQPZGro22RYPA3kAtjfWOY0Ij+MoQtZ7MeB2w69zVfzlHWOAso31zCMI5URNxeR1hPJMWYdItKMzPsby7tl9sSiUvBZQGZhQN3NrQ1GFtp799gIS6JqYwmi4mtjvA+lXVU+Eesbi0XQQkKtK7VMk7HQpdJZdXh4OXnzjpDpkFMc+9C+caiyZ2nuXSajJrtYxuRJ6CjykwzoiOaYQXfserOaC/2D6ABrbjgTBBOa0uMxMvN40QS5pBuc9XFwFkObkeaVZXqaf0WJnEQ4YGI7aVexlIrXg1CSp1S346jkNYylI6jkFQIXmUs4NCGzeRhhfwCGofva7XTLBN3IAPNzPYLa5hM3MbKZiaZSRc0Sxd8Gk0VKbUMEvE6TS0QJEntRfXiGyWhI7HIAYem3sU+nkcqqSSlfx4MDy/rDJP+KzZXsFZMzgLsYVN/vkoVo0LOBpkuXtROGcmCRSJ3xaGunvm016tDZiDGSOx4R7BtdCGZOjF5wSRifgdiuMSKMevBFsb9rvEumhEIRlCbSZai0dF4O7Sd5XqrSBH2LlCIA+o0zCmVkFZhyMgHnhODIX+SAaqv90AXc7Y//ffOb569efH778dM3X37w33/86tU/f3bz7EP4uP3r324++aXapgKjCZ76RCo0ae3sWI0t8LOUhIRxTdkde7HLJe3kRDbF6Rjpu50ctLu99uN+S+aTlJiIq5JiimvjIoyz4VhtQi4kj9RMPEFBnEmSn/S4wM2VKEdc4y5JEZ3n3yWobPH65eevX/75TmltsKzXsh892q7a9be2YnAhtldwaSvWL9QAAd6qWtvVej0HXsKgkSoAqZ16ntOwrUdWHW7daS73UkdeB1pm6Zkm+MxcklOktmqPUhfIjxOmdPwJhEjKVqGQG2gX/fInb/7ytIDS1tWU4iygnGxjxGJmMpBTH4wmxnh4NGp3jUl3LPwk4+VwavJvAOvV03gqOFUQxNj/jgHFHYqZUYEbUcUITw2C3XNjAZGYi8OCs62jqoXuOeWylo+jMJ6rE8kJZNQ1yDSbrJVYk2skT9pzJelkrDVwpevN379YJ6A30ska2vQJJE9DBY4i5fg6aj/716uvPnr14sXt85/ffPWjwgx6nfRlAPwcoik/TN0eUp6uN0sUdDyWxrywvp+rL5ooLlKHIQ04c+wH8g6lRwhEbTGb/I16Myh5cmJpL6AvUdC7mHWvuQ5s5wCZRQLoOcdQbcPs7pkSVA7P7Pqf3/7u9tdf3H76/M1Hf7r5+I83n3z6+uXv3zz7g4q626fPb3/6TGf55UIgdRE3WtUEGvKJxDVENBqidhtfv/8LIwi5AS2HEcuM9PX7n+UmE4rK5iSbGVq6VJGiCivQvLCQM3KqpDoU5FIR1QC0RQlrpAhdxMI5dbNF3quKqUTcSca3epNqzIgRQjf1bdhJEZwJbyqnRVRJPXxo1e265iptxBammOVM/9gPp9BkJAz5ALEEKUh9s0CGles97g9bzX4GUUoMI49Ewg3VB0paSlFSeiwZJa6WowAXGj839sWb0QpslZXOnEtjpn58OW16Iv2tf6ooIFA7jiLVEAX6OX8cz6EZTl7p7ubLl8tc/3ugetV8R5yNe50iH8Y5LhTCIlsQJF+mJs1SaarHxDuPamcPhGmyIfAKr51gDv1irzqtS+grI1PknW4UhdHa5JNxEtgAOmnIKGZm8RQjz1R13V52VgkhSXjph7r56R2GHeITXu4528ykB+HlvWXh7MuK9tjQ97Qxy109UrNkE+Tf9IWj/L9P+srZmlEEjZV4ACz9Bp9cXEdw3y2pjdqKFBQ3QFhdP5Lv04jxJyIT6C9FOUkpJ6pDfSIuXOpDjk+cxrYiAMDMz24W1ExCl6u/rIR+n85oyWuhlcR3cRKw5XyuWrheOU8RpeWAXEODmZsBkuL0B1A21DtKmdmUw0IuTeXF2yWjZ+e8rGIPp5h4ZGpV3SmpV7c8a6f6iJBG1bbhf+zW65a1LV4+9eSQOSi5KrmImRxY9hdP53+WdoZdKh0AAA=="

$c = $ws.Range("A1").Comment
$c.Text($newCommentText)

# 3) Update number formatting used by the yearly observation values
#    (numFmtId 166, applied to B27:B36) from "0.000" to "###0.000".
$ws.Range("B27:B36").NumberFormat = "###0.000"

# 4) Cell text + value corrections
$ws.Range("A11").Value = "Function Information"
$ws.Range("B21").Value = 0.2499825759175085
